# "hide filter in locations lists"
# Update the ImportTicketVisitSalesTemplate sample sheet:
#  - the "start date" / "end date" column headers gain an explicit
#    date-format hint for users filling in the template
#  - the sheet default font moves from Arial to Calibri
#  - the stale C2 selection left over from editing is cleared back to A1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header text for the start/end date columns (M1 / N1) so the
# expected format is explicit to whoever fills in the template.
$ws.Range("M1").Value = "start date (YYYY-MM-DD)"
$ws.Range("N1").Value = "end date (YYYY-MM-DD) "

# Switch the workbook's default (Normal) font from Arial to Calibri.
$wb.Styles.Item("Normal").Font.Name = "Calibri"

# Re-assert the yyyy-mm-dd date format on the date columns: changing the
# Normal style's font can otherwise clobber any explicit number format
# already applied to those cells.
$ws.Range("M1:N1").NumberFormat = "yyyy\-mm\-dd;@"

# Let the start/end date columns grow to fit their new, longer headers.
$ws.Columns.Item(13).AutoFit()
$ws.Columns.Item(14).AutoFit()

# Reset the lingering cell selection back to A1 (the default).
$ws.Range("A1").Select() | Out-Null
